# This script re-applies a reordering of the data rows (2-82) in the
# "Avverkningsanmälningar" sheet together with a global bump of the
# "Förändrad" date (column C) from serial 46062 to 46063.
#
# The underlying edit is: every data row keeps all of its own field
# values (A..Z), but a number of rows are relocated to a different row
# position (this corresponds to the list being freshly re-sorted /
# re-generated upstream). Column C (the "last changed" date) is set to
# the same new value (46063) for every data row, regardless of any
# relocation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 82
$rowCount = $lastRow - $firstRow + 1

# Destination row (array index, 0-based, 0 => row 2) -> source row number
# in the *current* (before) sheet that should end up there.
$sourceRowForDest = @(
    2,3,4,5,7,6,8,11,9,10,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,
    38,39,52,58,68,79,67,74,77,36,35,40,42,46,45,72,69,31,56,32,44,48,
    71,50,49,60,70,78,30,53,54,64,57,55,73,29,41,62,47,59,61,27,28,43,
    63,65,75,37,51,66,76,80,81,82,34,33
)

# --- Read the current state of the two column blocks ---------------------
# Columns A..R (1..18) hold plain values (text/number/date).
# Columns S..Z (19..26) hold HYPERLINK() formulas.
$valuesRange = $ws.Range("A$firstRow" + ":R$lastRow")
$formulaRange = $ws.Range("S$firstRow" + ":Z$lastRow")

$oldValues = $valuesRange.Value2
$oldFormulas = $formulaRange.Formula

# --- Build the reordered arrays -------------------------------------------
$newValues = New-Object 'object[,]' $rowCount,18
$newFormulas = New-Object 'object[,]' $rowCount,8

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $sourceRowForDest[$i]
    $srcIdx = $srcRow - $firstRow + 1   # 1-based index into $oldValues/$oldFormulas

    for ($c = 1; $c -le 18; $c++) {
        $newValues[$i, $c - 1] = $oldValues[$srcIdx, $c]
    }
    for ($c = 1; $c -le 8; $c++) {
        $newFormulas[$i, $c - 1] = $oldFormulas[$srcIdx, $c]
    }

    # Column C is the 3rd column of the A..R block (index 2, 0-based)
    $newValues[$i, 2] = 46063
}

# --- Write the reordered data back ----------------------------------------
$valuesRange.Value2 = $newValues
$formulaRange.Formula = $newFormulas
